# Updates cryptos list data (Coin name / Link / Price / Volume(1h)) to match
# the latest scrape, as produced by the "Updated cryptos list" GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.250.82'
$ws.Range("E2").Value = '  -0.92%  '

# Row 3
$ws.Range("D3").Value = '3.279.34'
$ws.Range("E3").Value = '  +0.27%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.36%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.57%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.605'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.13%  '

# Row 9
$ws.Range("E9").Value = '  -2.05%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.08%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.411'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.29%  '

# Row 12
$ws.Range("D12").Value = '3.856.03'
$ws.Range("E12").Value = '  +0.52%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.138'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.61%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.47%  '

# Row 15
$ws.Range("D15").Value = '68.206.34'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000169'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.55%  '

# Row 17
$ws.Range("D17").Value = '3.256.90'
$ws.Range("E17").Value = '  -1.21%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.06%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.75%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '416.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.48%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.34%  '

# Row 22
$ws.Range("E22").Value = '  +0.09%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.81%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.511'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.13%  '

# Row 25
$ws.Range("E25").Value = '  -1.50%  '

# Row 26
$ws.Range("E26").Value = '  -1.31%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.82%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.08%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.53%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.79'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.92%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.92%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.93'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.90%  '

# Row 33
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.65%  '

# Row 34
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.08%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '164.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.17%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.89%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.14%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.28'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.41%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.802'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.40%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.51'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.88%  '

# Row 41
$ws.Range("E41").Value = '  -4.14%  '

# Row 42
$ws.Range("D42").Value = '2.670.94'
$ws.Range("E42").Value = '  +2.35%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.23%  '

# Row 44
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.38%  '

# Row 45
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0680'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.50%  '

# Row 46
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '337.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.45%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.74'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.13%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0275'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.80%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.06%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.101'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.26%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.975'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.95%  '
